$d = $word.ActiveDocument
$d.Content.Font.ItalicBi = $true
